{"js": "// Two changes in the introductory paragraph (5th paragraph of the body):\n//   1) The whole paragraph's font size is set to 14pt (w:sz 28 half-points).\n//   2) The trailing phrase \"paragraph to \" is replaced with \"paraset \".\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph index 4 is \"You also can create a structure array ... is an\n// introductory paragraph to \" (the only non-empty paragraph in the doc).\nconst target = paragraphs.items[4];\ntarget.font.size = 14;\nawait context.sync();\n\nconst results = body.search(\"paragraph to \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"paraset \", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Two changes in the introductory paragraph (5th paragraph of the body):\n#   1) The whole paragraph's font size is set to 14pt (w:sz 28 half-points).\n#   2) The trailing phrase \"paragraph to \" is replaced with \"paraset \".\n$d = $word.ActiveDocument\n\n# Paragraph 5 (1-based) is \"You also can create a structure array ... is an\n# introductory paragraph to \" (the only non-empty paragraph in the doc).\n$p = $d.Paragraphs.Item(5)\n$p.Range.Font.Size = 14\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"paragraph to \"\n$find.Replacement.Text = \"paraset \"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
